# edit.ps1
# -----------------------------------------------------------------------
# 1) Slide 5 has a 3-column table whose table style is switched from the
#    deck's custom "Table_0" style to a built-in PowerPoint gallery style
#    (same GUID PowerPoint itself writes when you pick a style from the
#    Table Styles gallery on the ribbon).
# 2) The presentation's theme colour scheme (currently the "Red Violet"
#    palette of the "Integral" theme, used by the slide master / all
#    slides) is switched to the standard "Office" colour palette, i.e.
#    the same color values PowerPoint would write after picking the
#    built-in "Office Theme" design / "Office" color scheme from the
#    Design tab.
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 5 -----------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{00E5F180-FE4F-4E7C-96B7-8091AC5CBE88}")

# --- 2) Re-colour the theme to the "Office" palette ---------------------
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme

$themeColors.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$themeColors.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$themeColors.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$themeColors.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$themeColors.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$themeColors.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$themeColors.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$themeColors.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$themeColors.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$themeColors.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$themeColors.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$themeColors.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
